$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 399, shifting old rows 399-496 down to 402-499
$ws.Rows("399:401").Insert()

# Row 399
$ws.Range("A399").Value = 7
$ws.Range("B399").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C399").Value = "Ñuble"
$ws.Range("D399").Value = 44543
$ws.Range("E399").Value = 16
$ws.Range("F399").Value = "Fruta"
$ws.Range("G399").Value = 100102
$ws.Range("H399").Value = "Cítricos"
$ws.Range("I399").Value = 100102003
$ws.Range("J399").Value = "Limón"
$ws.Range("K399").Value = "Sin especificar"
$ws.Range("L399").Value = "1a amarillo"
$ws.Range("M399").Value = 160
$ws.Range("N399").Value = 8000
$ws.Range("O399").Value = 8500
$ws.Range("P399").Value = 8250
$ws.Range("Q399").Value = "`$/malla 16 kilos"
$ws.Range("R399").Value = "Provincia de Limarí"
$ws.Range("S399").Value = 516
$ws.Range("T399").Value = 16

# Row 400
$ws.Range("A400").Value = 7
$ws.Range("B400").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C400").Value = "Ñuble"
$ws.Range("D400").Value = 44543
$ws.Range("E400").Value = 16
$ws.Range("F400").Value = "Fruta"
$ws.Range("G400").Value = 100102
$ws.Range("H400").Value = "Cítricos"
$ws.Range("I400").Value = 100102003
$ws.Range("J400").Value = "Limón"
$ws.Range("K400").Value = "Sin especificar"
$ws.Range("L400").Value = "1a plateado"
$ws.Range("M400").Value = 160
$ws.Range("N400").Value = 10000
$ws.Range("O400").Value = 11000
$ws.Range("P400").Value = 10500
$ws.Range("Q400").Value = "`$/malla 16 kilos"
$ws.Range("R400").Value = "Región de O'Higgins"
$ws.Range("S400").Value = 656
$ws.Range("T400").Value = 16

# Row 401
$ws.Range("A401").Value = 7
$ws.Range("B401").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C401").Value = "Ñuble"
$ws.Range("D401").Value = 44543
$ws.Range("E401").Value = 16
$ws.Range("F401").Value = "Fruta"
$ws.Range("G401").Value = 100102
$ws.Range("H401").Value = "Cítricos"
$ws.Range("I401").Value = 100102003
$ws.Range("J401").Value = "Limón"
$ws.Range("K401").Value = "Sin especificar"
$ws.Range("L401").Value = "2a amarillo"
$ws.Range("M401").Value = 160
$ws.Range("N401").Value = 7000
$ws.Range("O401").Value = 7500
$ws.Range("P401").Value = 7250
$ws.Range("Q401").Value = "`$/malla 16 kilos"
$ws.Range("R401").Value = "Provincia de Limarí"
$ws.Range("S401").Value = 453
$ws.Range("T401").Value = 16
